$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nlgn2"
$ws.Range("C2").Value = "Nrxn1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.258986
$ws.Range("H2").Value = 9.776958
$ws.Range("I2").Value = 0.0928452675546778
$ws.Range("J2").Value = 0.09284526755467781
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2191816666666667
$ws.Range("N2").Value = 0.657545
$ws.Range("O2").Value = 0.8459226744507667
$ws.Range("P2").Value = 0.8459226744507669
$ws.Range("Q2").Value = 0.7143099831233335
$ws.Range("R2").Value = 6.428789848110001
$ws.Range("S2").Value = 0.07853991703995004
$ws.Range("T2").Value = 0.07853991703995007

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nlgn2"
$ws.Range("C3").Value = "Nrxn1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.258986
$ws.Range("H3").Value = 9.776958
$ws.Range("I3").Value = 0.0928452675546778
$ws.Range("J3").Value = 0.09284526755467781
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03648100000000001
$ws.Range("N3").Value = 0.109443
$ws.Range("O3").Value = 0.1407969268413801
$ws.Range("P3").Value = 0.1407969268413801
$ws.Range("Q3").Value = 0.118891068266
$ws.Range("R3").Value = 1.070019614394
$ws.Range("S3").Value = 0.01307232834346433
$ws.Range("T3").Value = 0.01307232834346434

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nlgn2"
$ws.Range("C4").Value = "Nrxn1"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.258986
$ws.Range("H4").Value = 9.776958
$ws.Range("I4").Value = 0.0928452675546778
$ws.Range("J4").Value = 0.09284526755467781
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.003441
$ws.Range("N4").Value = 0.010323
$ws.Range("O4").Value = 0.0132803987078531
$ws.Range("P4").Value = 0.0132803987078531
$ws.Range("Q4").Value = 0.011214170826
$ws.Range("R4").Value = 0.100927537434
$ws.Range("S4").Value = 0.001233022171263418
$ws.Range("T4").Value = 0.001233022171263418

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nlgn2"
$ws.Range("C5").Value = "Nrxn1"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.80139666666667
$ws.Range("H5").Value = 62.40418999999999
$ws.Range("I5").Value = 0.5926110879358332
$ws.Range("J5").Value = 0.5926110879358333
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2191816666666667
$ws.Range("N5").Value = 0.657545
$ws.Range("O5").Value = 0.8459226744507667
$ws.Range("P5").Value = 0.8459226744507669
$ws.Range("Q5").Value = 4.559284790394445
$ws.Range("R5").Value = 41.03356311355
$ws.Range("S5").Value = 0.5013031564158585
$ws.Range("T5").Value = 0.5013031564158587

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Nlgn2"
$ws.Range("C6").Value = "Nrxn1"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 20.80139666666667
$ws.Range("H6").Value = 62.40418999999999
$ws.Range("I6").Value = 0.5926110879358332
$ws.Range("J6").Value = 0.5926110879358333
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03648100000000001
$ws.Range("N6").Value = 0.109443
$ws.Range("O6").Value = 0.1407969268413801
$ws.Range("P6").Value = 0.1407969268413801
$ws.Range("Q6").Value = 0.7588557517966668
$ws.Range("R6").Value = 6.82970176617
$ws.Range("S6").Value = 0.08343781999349217
$ws.Range("T6").Value = 0.0834378199934922

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Nlgn2"
$ws.Range("C7").Value = "Nrxn1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 20.80139666666667
$ws.Range("H7").Value = 62.40418999999999
$ws.Range("I7").Value = 0.5926110879358332
$ws.Range("J7").Value = 0.5926110879358333
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.003441
$ws.Range("N7").Value = 0.010323
$ws.Range("O7").Value = 0.0132803987078531
$ws.Range("P7").Value = 0.0132803987078531
$ws.Range("Q7").Value = 0.07157760593
$ws.Range("R7").Value = 0.6441984533699999
$ws.Range("S7").Value = 0.007870111526482456
$ws.Range("T7").Value = 0.007870111526482459

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Nlgn2"
$ws.Range("C8").Value = "Nrxn1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.70406433333333
$ws.Range("H8").Value = 32.112193
$ws.Range("I8").Value = 0.3049481393755043
$ws.Range("J8").Value = 0.3049481393755044
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2191816666666667
$ws.Range("N8").Value = 0.657545
$ws.Range("O8").Value = 0.8459226744507667
$ws.Range("P8").Value = 0.8459226744507669
$ws.Range("Q8").Value = 2.346134660687223
$ws.Range("R8").Value = 21.11521194618501
$ws.Range("S8").Value = 0.2579625456293118
$ws.Range("T8").Value = 0.2579625456293119

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Nlgn2"
$ws.Range("C9").Value = "Nrxn1"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.70406433333333
$ws.Range("H9").Value = 32.112193
$ws.Range("I9").Value = 0.3049481393755043
$ws.Range("J9").Value = 0.3049481393755044
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.03648100000000001
$ws.Range("N9").Value = 0.109443
$ws.Range("O9").Value = 0.1407969268413801
$ws.Range("P9").Value = 0.1407969268413801
$ws.Range("Q9").Value = 0.3904949709443334
$ws.Range("R9").Value = 3.514454738499001
$ws.Range("S9").Value = 0.04293576087006786
$ws.Range("T9").Value = 0.04293576087006788

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Nlgn2"
$ws.Range("C10").Value = "Nrxn1"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.70406433333333
$ws.Range("H10").Value = 32.112193
$ws.Range("I10").Value = 0.3049481393755043
$ws.Range("J10").Value = 0.3049481393755044
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.003441
$ws.Range("N10").Value = 0.010323
$ws.Range("O10").Value = 0.0132803987078531
$ws.Range("P10").Value = 0.0132803987078531
$ws.Range("Q10").Value = 0.03683268537100001
$ws.Range("R10").Value = 0.3314941683390001
$ws.Range("S10").Value = 0.004049832876124653
$ws.Range("T10").Value = 0.004049832876124655

$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Nlgn2"
$ws.Range("C11").Value = "Nrxn1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3368143333333333
$ws.Range("H11").Value = 1.010443
$ws.Range("I11").Value = 0.009595505133984546
$ws.Range("J11").Value = 0.00959550513398455
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2191816666666667
$ws.Range("N11").Value = 0.657545
$ws.Range("O11").Value = 0.8459226744507667
$ws.Range("P11").Value = 0.8459226744507669
$ws.Range("Q11").Value = 0.07382352693722223
$ws.Range("R11").Value = 0.664411742435
$ws.Range("S11").Value = 0.00811705536564627
$ws.Range("T11").Value = 0.008117055365646275

$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Nlgn2"
$ws.Range("C12").Value = "Nrxn1"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.3368143333333333
$ws.Range("H12").Value = 1.010443
$ws.Range("I12").Value = 0.009595505133984546
$ws.Range("J12").Value = 0.00959550513398455
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.03648100000000001
$ws.Range("N12").Value = 0.109443
$ws.Range("O12").Value = 0.1407969268413801
$ws.Range("P12").Value = 0.1407969268413801
$ws.Range("Q12").Value = 0.01228732369433333
$ws.Range("R12").Value = 0.110585913249
$ws.Range("S12").Value = 0.001351017634355709
$ws.Range("T12").Value = 0.00135101763435571

$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Nlgn2"
$ws.Range("C13").Value = "Nrxn1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.3368143333333333
$ws.Range("H13").Value = 1.010443
$ws.Range("I13").Value = 0.009595505133984546
$ws.Range("J13").Value = 0.00959550513398455
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.003441
$ws.Range("N13").Value = 0.010323
$ws.Range("O13").Value = 0.0132803987078531
$ws.Range("P13").Value = 0.0132803987078531
$ws.Range("Q13").Value = 0.001158978121
$ws.Range("R13").Value = 0.010430803089
$ws.Range("S13").Value = 0.0001274321339825661
$ws.Range("T13").Value = 0.0001274321339825662
